# update logic import store new
#
# Corrects the PRODUCT QUANTITY (column I) values on the "Products" sheet
# that were produced by the previous (buggy) import logic, and leaves the
# sheet scrolled/selected where the user ended up after making the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

# --- Fix PRODUCT QUANTITY values (column I) per corrected import logic ---
$ws.Range("I9").Value  = 0
$ws.Range("I26").Value = 0
$ws.Range("I27").Value = 30
$ws.Range("I28").Value = 20
$ws.Range("I39").Value = 30
$ws.Range("I45").Value = 0
$ws.Range("I49").Value = 0

# --- Update the view: scroll so column D is left-most, then select I53 ---
$window = $excel.ActiveWindow
$window.ScrollColumn = 4   # column D becomes the left-most visible column
$window.ScrollRow = 43     # row 43 becomes the top-most visible row

$ws.Range("I53").Select()
